# Updated symbol list on Wed Dec 21 07:09:30 UTC 2022 with GitHub Actions
#
# This script applies the 2022-12-21 07:09 UTC data refresh to the "cryptos"
# worksheet:
#   - Column G (Hora / hour-of-update) moves from "6" to "7" for every data
#     row (rows 2-51).
#   - Column D (Price) is refreshed with newer quotes for the rows whose
#     price changed.
#   - Rows 10-19 (coin list) rotate: "One" (previously row 19) becomes the
#     new row 10, and WazirX ... UpBots each shift down by one row, with
#     their Link/Price/Volume columns following them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D (Price) new values -------------------------------------------------
# These look numeric, so we force the cell's number format to Text ("@") first;
# otherwise Excel would silently convert the inline string into a numeric value.
$dValues = [ordered]@{
    "D2"  = "248.38"
    "D3"  = "22.67"
    "D4"  = "5.408"
    "D5"  = "0.05703"
    "D6"  = "3.407"
    "D7"  = "6.315"
    "D8"  = "0.8128"
    "D9"  = "0.9226"
    "D10" = "0.01125"
    "D11" = "0.1413"
    "D12" = "0.07445"
    "D13" = "0.03077"
    "D14" = "0.03018"
    "D15" = "0.09380"
    "D16" = "3.750"
    "D17" = "0.001578"
    "D18" = "0.04769"
    "D19" = "0.01828"
    "D20" = "0.006463"
    "D21" = "0.004998"
    "D22" = "0.001026"
    "D23" = "0.0001501"
    "D27" = "0.1298"
    "D40" = "0.03991"
    "D41" = "0.006853"
    "D42" = "0.1067"
    "D44" = "0.007509"
    "D45" = "0.00005802"
    "D48" = "0.2115"
}

# --- Column B / C / E new text values (row rotation for rows 10-19) --------------
$textValues = [ordered]@{
    "B10" = "One"
    "C10" = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
    "E10" = "9OneONE"
    "B11" = "WazirX"
    "C11" = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
    "E11" = "10WazirXWRX"
    "B12" = "MandalaExchangeToken"
    "C12" = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
    "E12" = "11MandalaExchangeTokenMDX"
    "B13" = "LiechtensteinCryptoassetsExchange"
    "C13" = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
    "E13" = "12LiechtensteinCryptoassetsExchangeLCX"
    "B14" = "BitrueCoin"
    "C14" = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
    "E14" = "13BitrueCoinBTR"
    "B15" = "BitMartToken"
    "C15" = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
    "E15" = "14BitMartTokenBMX"
    "B16" = "MCDex"
    "C16" = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
    "E16" = "15MCDexMCB"
    "B17" = "BitForexToken"
    "C17" = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
    "E17" = "16BitForexTokenBF"
    "B18" = "CoinExToken"
    "C18" = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
    "E18" = "17CoinExTokenCET"
    "B19" = "UpBots"
    "C19" = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
    "E19" = "18UpBotsUBXTBestin24h"
}

# Force the whole Price column (data rows) to Text format so that the new
# numeric-looking inline-string values are not reinterpreted as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

foreach ($key in $dValues.Keys) {
    $ws.Range($key).Value = $dValues[$key]
}

foreach ($key in $textValues.Keys) {
    $ws.Range($key).Value = $textValues[$key]
}

# --- Column G (Hora) : every data row goes from 6 to 7 ---------------------------
$ws.Range("G2:G51").NumberFormat = "@"
for ($row = 2; $row -le 51; $row++) {
    $ws.Range("G$row").Value = "7"
}
